# add chapter 4 and 5: 装修 (renovation) and 楼层 (floor) columns

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Copy the header style (s="1") from J1 onto the two new header cells K1:L1
$ws.Range("J1").Copy() | Out-Null
$ws.Range("K1:L1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# Column K - 装修 (renovation level)
$ws.Cells.Item(1, 11).Value = "装修"
$ws.Cells.Item(2, 11).Value = "豪华装修"
$ws.Cells.Item(3, 11).Value = "精装修"
$ws.Cells.Item(4, 11).Value = "精装修"
$ws.Cells.Item(5, 11).Value = "简装修"
$ws.Cells.Item(6, 11).Value = "中装修"
$ws.Cells.Item(7, 11).Value = "精装修"
$ws.Cells.Item(8, 11).Value = "精装修"
$ws.Cells.Item(9, 11).Value = "简装修"
$ws.Cells.Item(10, 11).Value = "简装修"
$ws.Cells.Item(11, 11).Value = "简装修"
$ws.Cells.Item(12, 11).Value = "豪华装修"
$ws.Cells.Item(13, 11).Value = "简装修"
$ws.Cells.Item(14, 11).Value = "简装修"
$ws.Cells.Item(15, 11).Value = "简装修"
$ws.Cells.Item(16, 11).Value = "简装修"
$ws.Cells.Item(17, 11).Value = "精装修"
$ws.Cells.Item(18, 11).Value = "精装修"
$ws.Cells.Item(19, 11).Value = "简装修"
$ws.Cells.Item(20, 11).Value = "简装修"
$ws.Cells.Item(21, 11).Value = "简装修"
$ws.Cells.Item(22, 11).Value = "精装修"
$ws.Cells.Item(23, 11).Value = "豪华装修"
$ws.Cells.Item(24, 11).Value = "精装修"
$ws.Cells.Item(25, 11).Value = "精装修"
$ws.Cells.Item(26, 11).Value = "精装修"
$ws.Cells.Item(27, 11).Value = "精装修"
$ws.Cells.Item(28, 11).Value = "精装修"
$ws.Cells.Item(29, 11).Value = "精装修"
$ws.Cells.Item(30, 11).Value = "精装修"

# Column L - 楼层 (floor level)
$ws.Cells.Item(1, 12).Value = "楼层"
$ws.Cells.Item(2, 12).Value = "低层"
$ws.Cells.Item(3, 12).Value = "中层"
$ws.Cells.Item(4, 12).Value = "高层"
$ws.Cells.Item(5, 12).Value = "中层"
$ws.Cells.Item(6, 12).Value = "低层"
$ws.Cells.Item(7, 12).Value = "中层"
$ws.Cells.Item(8, 12).Value = "中层"
$ws.Cells.Item(9, 12).Value = "低层"
$ws.Cells.Item(10, 12).Value = "中层"
$ws.Cells.Item(11, 12).Value = "中层"
$ws.Cells.Item(12, 12).Value = "高层"
$ws.Cells.Item(13, 12).Value = "高层"
$ws.Cells.Item(14, 12).Value = "高层"
$ws.Cells.Item(15, 12).Value = "低层"
$ws.Cells.Item(16, 12).Value = "中层"
$ws.Cells.Item(17, 12).Value = "高层"
$ws.Cells.Item(18, 12).Value = "高层"
$ws.Cells.Item(19, 12).Value = "中层"
$ws.Cells.Item(20, 12).Value = "中层"
$ws.Cells.Item(21, 12).Value = "高层"
$ws.Cells.Item(22, 12).Value = "中层"
$ws.Cells.Item(23, 12).Value = "中层"
$ws.Cells.Item(24, 12).Value = "低层"
$ws.Cells.Item(25, 12).Value = "低层"
$ws.Cells.Item(26, 12).Value = "中层"
$ws.Cells.Item(27, 12).Value = "高层"
$ws.Cells.Item(28, 12).Value = "中层"
$ws.Cells.Item(29, 12).Value = "低层"
$ws.Cells.Item(30, 12).Value = "中层"

# Update the view: scroll back to top-left (drop topLeftCell="A19") and move
# the active selection to E12 instead of K26
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E12").Select() | Out-Null
